$d = $word.ActiveDocument

# Hybrid bold + color (#2C3E50) highlighting for quantitative impact metrics.
# Word COM Font.Color takes a BGR-packed long; 0x2C3E50 (RGB) -> 0x503E2C (BGR) = 5258796.
$metricColor = 5258796

function Set-MetricFormat($rng) {
    $rng.Font.Bold = 1
    $rng.Font.Color = $metricColor
}

# Within the paragraph identified by $anchorText (a unique substring of the
# paragraph's plain text), bold+color each substring in $metrics, in order,
# searching left-to-right so repeated values (e.g. "73.5%" appearing twice
# across different paragraphs) are only matched inside this paragraph.
function Highlight-Metrics($anchorText, $metrics) {
    $count = $d.Paragraphs.Count
    $target = $null
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($target -eq $null -and $p.Range.Text.Contains($anchorText)) {
            $target = $p
        }
    }
    if ($target -eq $null) {
        Write-Output "PARAGRAPH NOT FOUND: $anchorText"
        return
    }

    $searchStart = $target.Range.Start
    $paraEnd = $target.Range.End
    foreach ($sub in $metrics) {
        $r = $d.Range($searchStart, $paraEnd)
        $found = $r.Find.Execute($sub, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            Set-MetricFormat $r
            $searchStart = $r.End
        } else {
            Write-Output "METRIC NOT FOUND: $sub (in paragraph containing '$anchorText')"
        }
    }
}

Highlight-Metrics "Discovered systematic race coding errors" @("23%", "64%")

Highlight-Metrics "Utilized advanced sampling methods" @("±4.2%", "±2.1%", "71%", "87%")

Highlight-Metrics "Trigonometric algorithm for boundary estimation" @("73.5%", "$4.7M")

Highlight-Metrics "Built real-time FEC analysis systems" @("$2")

Highlight-Metrics "Modernized legacy ETL processes" @("57%")

Highlight-Metrics "Algorithmic innovation: Pioneered trigonometric boundary estimation" @("73.5%")

Highlight-Metrics "$4.7M savings enabled nonprofit access" @("$4.7M")

Highlight-Metrics "Platform impact: Built redistricting system serving" @("12,847")

Write-Output "DONE"
